$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = "CK203"
$ws.Range("B7").Value = 7939
$ws.Range("C7").Value = 43713
$ws.Range("C7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = "Nitinat Lake"
$ws.Range("E7").Value = 48.68626
$ws.Range("F7").Value = -124.836131
$ws.Range("G7").Value = "hatchery"
$ws.Range("H7").Value = "Nitinat Hatchery"
$ws.Range("I7").Value = "Caroline Sherry"
$ws.Range("J7").Value = 12507453321
$ws.Range("K7").Value = "NA"
$ws.Range("L7").Value = "PO Box 369, Port Alberni BC"

# Row 8
$ws.Range("A8").Value = "CK65"
$ws.Range("B8").Value = 7715
$ws.Range("C8").Value = 43697
$ws.Range("C8").NumberFormat = "d-mmm"
$ws.Range("D8").Value = "Wya Point"
$ws.Range("E8").Value = 48.970067
$ws.Range("F8").Value = -125.620663
$ws.Range("G8").Value = "sport"
$ws.Range("H8").Value = "Nitinat Hatchery"
$ws.Range("I8").Value = "Caroline Sherry"
$ws.Range("J8").Value = 12507453321
$ws.Range("K8").Value = "NA"
$ws.Range("L8").Value = "PO Box 369, Port Alberni BC"
$ws.Range("M8").Value = "captured by acquiantance; may not be returned"

# Row 9
$ws.Range("A9").Value = "CK262"
$ws.Range("B9").Value = 8207
$ws.Range("C9").Value = "NA"
$ws.Range("D9").Value = "Columbia River"
$ws.Range("G9").Value = "sport"
$ws.Range("M9").Value = "released; contacted indirectly by Joe Smith (NOAA)"

# Row 10
$ws.Range("A10").Value = "CK01"
$ws.Range("B10").Value = 7689
$ws.Range("C10").Value = 43713
$ws.Range("C10").NumberFormat = "d-mmm"
$ws.Range("D10").Value = "Rock Creek (Columbia River)"
$ws.Range("E10").Value = 45.699452
$ws.Range("F10").Value = -120.379026
$ws.Range("G10").Value = "gillnet"
$ws.Range("H10").Value = "WDFW"
$ws.Range("I10").Value = "Rick Heitz"
$ws.Range("J10").Value = 13606095688
$ws.Range("K10").Value = "NA"
$ws.Range("L10").Value = "5525 S 11th St., Ridgefield WA, 98642"
$ws.Range("M10").Value = "captured by unknown FN group"

# Update selection to match target (A11 selected after data entry)
$ws.Range("A11").Select() | Out-Null
